$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "SmartCash ist sehr neu in der Szene, so dass die Beschaffung vielleicht ein bisschen komplizierter ist, als einfach ein Kauf bei Coinbase. Aber keine Sorge, wir sind der richtige Partner für Sie.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "SmartCash ist sehr neu in der Szene, so dass es ein wenig komplexer sein kann, als einfach nur von Coinbase zu kaufen, aber wir helfen Dir.",
    2)

$d.Content.Find.Execute(
    "VERWENDUNG EINES MINERS:",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "BENUTZ DEN MINER:",
    2)

$d.Content.Find.Execute(
    "Für Support hinsichtlich Mining bitte bei ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Für Mining Support bitte hier anmelden ",
    2)

$d.Content.Find.Execute(
    " melden",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "",
    2)
